$wb = $excel.ActiveWorkbook

# --- Sheet 1: Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3,4).Value = 22
$ws1.Cells.Item(3,6).Value = 17
$ws1.Cells.Item(3,7).Value = 43.59
$ws1.Cells.Item(3,8).Value = 8.300000000000001
$ws1.Cells.Item(4,4).Value = 20
$ws1.Cells.Item(4,6).Value = 15
$ws1.Cells.Item(4,7).Value = 42.86
$ws1.Cells.Item(4,8).Value = 8.1
$ws1.Cells.Item(5,4).Value = 19
$ws1.Cells.Item(5,6).Value = 13
$ws1.Cells.Item(5,7).Value = 40.63
$ws1.Cells.Item(5,8).Value = 7.2

# --- Sheet 2: Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3,5).Value = 17
$ws2.Cells.Item(4,5).Value = 15
$ws2.Cells.Item(5,5).Value = 13

# --- Sheet 3: Estadisticos Final ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3,4).Value = 22
$ws3.Cells.Item(3,6).Value = 17
$ws3.Cells.Item(3,7).Value = 43.59
$ws3.Cells.Item(3,8).Value = 8.300000000000001
$ws3.Cells.Item(4,4).Value = 20
$ws3.Cells.Item(4,6).Value = 15
$ws3.Cells.Item(4,7).Value = 42.86
$ws3.Cells.Item(4,8).Value = 8.1
$ws3.Cells.Item(5,4).Value = 19
$ws3.Cells.Item(5,6).Value = 13
$ws3.Cells.Item(5,7).Value = 40.63
$ws3.Cells.Item(5,8).Value = 7.2

# --- Sheet 4: Rescatables ---
$ws4 = $wb.Worksheets.Item(4)

# Remove obsolete rows 30-35 (old list had 35 rows, new list has 29)
$ws4.Range("A30:G35").Delete()

# Overwrite rows 2-29 with the new roster
$ws4.Cells.Item(2,1).Value = 20330051920223
$ws4.Cells.Item(2,2).Value = 'CASTELLANOS'
$ws4.Cells.Item(2,3).Value = 'TEQUIHUATLE'
$ws4.Cells.Item(2,4).Value = 'JENNIFER'
$ws4.Cells.Item(2,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(2,6).Value = '3ALCM'
$ws4.Cells.Item(2,7).Value = 6

$ws4.Cells.Item(3,1).Value = 20330051920245
$ws4.Cells.Item(3,2).Value = 'PONCE'
$ws4.Cells.Item(3,3).Value = 'GOMEZ'
$ws4.Cells.Item(3,4).Value = 'ALETHIA LUCIA'
$ws4.Cells.Item(3,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(3,6).Value = '3ALCM'
$ws4.Cells.Item(3,7).Value = 6

$ws4.Cells.Item(4,1).Value = 20330051920250
$ws4.Cells.Item(4,2).Value = 'ROMAN'
$ws4.Cells.Item(4,3).Value = 'ANTONIO'
$ws4.Cells.Item(4,4).Value = 'FABIOLA'
$ws4.Cells.Item(4,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(4,6).Value = '3ALCM'
$ws4.Cells.Item(4,7).Value = 6

$ws4.Cells.Item(5,1).Value = 20330051920252
$ws4.Cells.Item(5,2).Value = 'SANCHEZ'
$ws4.Cells.Item(5,3).Value = 'PEREZ'
$ws4.Cells.Item(5,4).Value = 'ARLET'
$ws4.Cells.Item(5,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(5,6).Value = '3ALCM'
$ws4.Cells.Item(5,7).Value = 6

$ws4.Cells.Item(6,1).Value = 20330051920253
$ws4.Cells.Item(6,2).Value = 'SANCHEZ'
$ws4.Cells.Item(6,3).Value = 'QUIAHUA'
$ws4.Cells.Item(6,4).Value = 'ROSARIO'
$ws4.Cells.Item(6,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(6,6).Value = '3ALCM'
$ws4.Cells.Item(6,7).Value = 6

$ws4.Cells.Item(7,1).Value = 20330051920156
$ws4.Cells.Item(7,2).Value = 'ANTONIO'
$ws4.Cells.Item(7,3).Value = 'GUERRA'
$ws4.Cells.Item(7,4).Value = 'LUIS YAEL'
$ws4.Cells.Item(7,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(7,6).Value = '3ARHM'
$ws4.Cells.Item(7,7).Value = 6

$ws4.Cells.Item(8,1).Value = 20330051920158
$ws4.Cells.Item(8,2).Value = 'BAUTISTA'
$ws4.Cells.Item(8,3).Value = 'DIAZ'
$ws4.Cells.Item(8,4).Value = 'DINA BERENICE'
$ws4.Cells.Item(8,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(8,6).Value = '3ARHM'
$ws4.Cells.Item(8,7).Value = 6

$ws4.Cells.Item(9,1).Value = 20330051920306
$ws4.Cells.Item(9,2).Value = 'PARRA'
$ws4.Cells.Item(9,3).Value = 'FLORES'
$ws4.Cells.Item(9,4).Value = 'SUEMI'
$ws4.Cells.Item(9,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(9,6).Value = '3ARHM'
$ws4.Cells.Item(9,7).Value = 6

$ws4.Cells.Item(10,1).Value = 20330051920178
$ws4.Cells.Item(10,2).Value = 'QUIRIZ'
$ws4.Cells.Item(10,3).Value = 'RAMOS'
$ws4.Cells.Item(10,4).Value = 'MONICA'
$ws4.Cells.Item(10,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(10,6).Value = '3ARHM'
$ws4.Cells.Item(10,7).Value = 6

$ws4.Cells.Item(11,1).Value = 20330051920042
$ws4.Cells.Item(11,2).Value = 'ANASTACIO'
$ws4.Cells.Item(11,3).Value = 'ROMERO'
$ws4.Cells.Item(11,4).Value = 'HIRAM FABIAN'
$ws4.Cells.Item(11,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(11,6).Value = '3BEM'
$ws4.Cells.Item(11,7).Value = 6

$ws4.Cells.Item(12,1).Value = 20330051920374
$ws4.Cells.Item(12,2).Value = 'BERISTAIN'
$ws4.Cells.Item(12,3).Value = 'APALE'
$ws4.Cells.Item(12,4).Value = 'JOSE ISAIAS'
$ws4.Cells.Item(12,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(12,6).Value = '3ALCM'
$ws4.Cells.Item(12,7).Value = 6

$ws4.Cells.Item(13,1).Value = 20330051920224
$ws4.Cells.Item(13,2).Value = 'COLOHUA'
$ws4.Cells.Item(13,3).Value = 'RAMIREZ'
$ws4.Cells.Item(13,4).Value = 'FERNANDA'
$ws4.Cells.Item(13,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(13,6).Value = '3ALCM'
$ws4.Cells.Item(13,7).Value = 6

$ws4.Cells.Item(14,1).Value = 20330051920375
$ws4.Cells.Item(14,2).Value = 'FLORES'
$ws4.Cells.Item(14,3).Value = 'IXMATLAHUA'
$ws4.Cells.Item(14,4).Value = 'JENIFER'
$ws4.Cells.Item(14,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(14,6).Value = '3ALCM'
$ws4.Cells.Item(14,7).Value = 6

$ws4.Cells.Item(15,1).Value = 20330051920228
$ws4.Cells.Item(15,2).Value = 'HERNANDEZ'
$ws4.Cells.Item(15,3).Value = 'GALEOTE'
$ws4.Cells.Item(15,4).Value = 'GERMAN ISAI'
$ws4.Cells.Item(15,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(15,6).Value = '3ALCM'
$ws4.Cells.Item(15,7).Value = 6

$ws4.Cells.Item(16,1).Value = 20330051920232
$ws4.Cells.Item(16,2).Value = 'JUSTO'
$ws4.Cells.Item(16,3).Value = 'LORENZO'
$ws4.Cells.Item(16,4).Value = 'EVELYN'
$ws4.Cells.Item(16,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(16,6).Value = '3ALCM'
$ws4.Cells.Item(16,7).Value = 6

$ws4.Cells.Item(17,1).Value = 20330051920241
$ws4.Cells.Item(17,2).Value = 'MENDOZA'
$ws4.Cells.Item(17,3).Value = 'LEON'
$ws4.Cells.Item(17,4).Value = 'ITZEL'
$ws4.Cells.Item(17,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(17,6).Value = '3ALCM'
$ws4.Cells.Item(17,7).Value = 6

$ws4.Cells.Item(18,1).Value = 20330051920257
$ws4.Cells.Item(18,2).Value = 'VARGAS'
$ws4.Cells.Item(18,3).Value = 'TETLA'
$ws4.Cells.Item(18,4).Value = 'SELINA'
$ws4.Cells.Item(18,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(18,6).Value = '3ALCM'
$ws4.Cells.Item(18,7).Value = 6

$ws4.Cells.Item(19,1).Value = 20330051920258
$ws4.Cells.Item(19,2).Value = 'XOTLANIHUA'
$ws4.Cells.Item(19,3).Value = 'XOTLANIHUA'
$ws4.Cells.Item(19,4).Value = 'ERIKA'
$ws4.Cells.Item(19,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(19,6).Value = '3ALCM'
$ws4.Cells.Item(19,7).Value = 6

$ws4.Cells.Item(20,1).Value = 20330051920161
$ws4.Cells.Item(20,2).Value = 'CANSECO'
$ws4.Cells.Item(20,3).Value = 'LEAL'
$ws4.Cells.Item(20,4).Value = 'ANGELA'
$ws4.Cells.Item(20,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(20,6).Value = '3ARHM'
$ws4.Cells.Item(20,7).Value = 6

$ws4.Cells.Item(21,1).Value = 20330051920372
$ws4.Cells.Item(21,2).Value = 'FLORES'
$ws4.Cells.Item(21,3).Value = 'CERON'
$ws4.Cells.Item(21,4).Value = 'MARIA YAZMIN'
$ws4.Cells.Item(21,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(21,6).Value = '3ARHM'
$ws4.Cells.Item(21,7).Value = 6

$ws4.Cells.Item(22,1).Value = 20330051920172
$ws4.Cells.Item(22,2).Value = 'MATA'
$ws4.Cells.Item(22,3).Value = 'CANSECO'
$ws4.Cells.Item(22,4).Value = 'CRISTIAN ARTURO'
$ws4.Cells.Item(22,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(22,6).Value = '3ARHM'
$ws4.Cells.Item(22,7).Value = 6

$ws4.Cells.Item(23,1).Value = 20330051920254
$ws4.Cells.Item(23,2).Value = 'SANCHEZ'
$ws4.Cells.Item(23,3).Value = 'ROMERO'
$ws4.Cells.Item(23,4).Value = 'BERENICE'
$ws4.Cells.Item(23,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(23,6).Value = '3ARHM'
$ws4.Cells.Item(23,7).Value = 6

$ws4.Cells.Item(24,1).Value = 20330051920044
$ws4.Cells.Item(24,2).Value = 'BERNARDO'
$ws4.Cells.Item(24,3).Value = 'CONCHOA'
$ws4.Cells.Item(24,4).Value = 'URIEL'
$ws4.Cells.Item(24,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(24,6).Value = '3BEM'
$ws4.Cells.Item(24,7).Value = 6

$ws4.Cells.Item(25,1).Value = 20330051920046
$ws4.Cells.Item(25,2).Value = 'CID'
$ws4.Cells.Item(25,3).Value = 'VALENCIA'
$ws4.Cells.Item(25,4).Value = 'JESUS'
$ws4.Cells.Item(25,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(25,6).Value = '3BEM'
$ws4.Cells.Item(25,7).Value = 6

$ws4.Cells.Item(26,1).Value = 20330051920051
$ws4.Cells.Item(26,2).Value = 'CRUZ'
$ws4.Cells.Item(26,3).Value = 'MARROQUIN'
$ws4.Cells.Item(26,4).Value = 'JESUS'
$ws4.Cells.Item(26,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(26,6).Value = '3BEM'
$ws4.Cells.Item(26,7).Value = 6

$ws4.Cells.Item(27,1).Value = 20330051920101
$ws4.Cells.Item(27,2).Value = 'ORTIZ'
$ws4.Cells.Item(27,3).Value = 'ROSETE'
$ws4.Cells.Item(27,4).Value = 'SERGIO MARIANO'
$ws4.Cells.Item(27,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(27,6).Value = '3BEM'
$ws4.Cells.Item(27,7).Value = 6

$ws4.Cells.Item(28,1).Value = 20330051920105
$ws4.Cells.Item(28,2).Value = 'RODRIGUEZ'
$ws4.Cells.Item(28,3).Value = 'HERNANDEZ'
$ws4.Cells.Item(28,4).Value = 'CESAR OMAR'
$ws4.Cells.Item(28,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(28,6).Value = '3BEM'
$ws4.Cells.Item(28,7).Value = 6

$ws4.Cells.Item(29,1).Value = 20330051920111
$ws4.Cells.Item(29,2).Value = 'VALENTE'
$ws4.Cells.Item(29,3).Value = 'GAMEZ'
$ws4.Cells.Item(29,4).Value = 'ABIUD'
$ws4.Cells.Item(29,5).Value = 'GEOMETRÍA ANALÍTICA'
$ws4.Cells.Item(29,6).Value = '3BEM'
$ws4.Cells.Item(29,7).Value = 6

